{"js": "// Apply a series of targeted text corrections (punctuation, grammar, wording)\n// to the body of the document. Each entry replaces a unique substring with\n// its corrected counterpart; the surrounding text stays untouched.\nconst replacements = [\n  [\"In order to run our program the user\", \"In order to run our program, the user\"],\n  [\" With in this file we begin by importing \", \" With in this file, we begin by importing \"],\n  [\n    \"made the game to display in full screen. After the window and background the settings,\",\n    \"made the game display in full screen. After the window and background, the settings,\",\n  ],\n  [\n    \"begin upon launching from the players terminal, however\",\n    \"begin upon launching from the player\\u2019s terminal, however\",\n  ],\n  [\" A player looses a ship if \", \" A player loses a ship if \"],\n  [\n    \"the high score is checked and all settings are reset\",\n    \"the high score is checked, and all settings are reset\",\n  ],\n  [\"This file although very short allows\", \"This file, although very short allows\"],\n  [\"As mentioned it also resets\", \"As mentioned, it also resets\"],\n  [\n    \"represents the ships location. The file has\",\n    \"represents the ships\\u2019 location. The file has\",\n  ],\n  [\n    \"one of which updates the ships position based on the players input,\",\n    \"one of which updates the ship\\u2019s position based on the players input,\",\n  ],\n  [\n    \"speed up subsequent imports of the same module \",\n    \"speed up subsequent imports of the same module. \",\n  ],\n  [\n    \"Here we will talk about adjusting the speeds and sizes of bullets that we tried to get a challenging but fun user experience. While testing, during various steps of the build , we changed\",\n    \"Here we talked about and tested adjusting the speeds and sizes of bullets, we tried to get a challenging but fun user experience. While testing, during various steps of the build, we changed\",\n  ],\n  // \"complete\" is the unique remaining word of the \"Aug 5 - 10\" bullet; matched\n  // on its own (rather than together with the preceding sentence) so the\n  // existing grammar-check markers around it stay intact.\n  [\n    \"complete\",\n    \"complete, still unable to reach other members of group B\",\n  ],\n  [\"By the end of the project it made sense\", \"By the end of the project, it made sense\"],\n  [\"or any other real world tasks,\", \"or any other real-world tasks,\"],\n  [\n    \"We have also discussed about finding and implementing\",\n    \"We have also discussed finding and implementing\",\n  ],\n  [\"maybe the bullets do no disappear\", \"maybe the bullets do not disappear\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply a series of targeted text corrections (punctuation, grammar, wording)\n# throughout the document body. Each pair is a unique \"Find\" string and its\n# corrected \"Replace\" string; everything else in the document is left as-is.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"In order to run our program the user\", \"In order to run our program, the user\"),\n    @(\" With in this file we begin by importing \", \" With in this file, we begin by importing \"),\n    @(\"made the game to display in full screen. After the window and background the settings,\",\n      \"made the game display in full screen. After the window and background, the settings,\"),\n    @(\"begin upon launching from the players terminal, however\",\n      \"begin upon launching from the player\" + [char]0x2019 + \"s terminal, however\"),\n    @(\" A player looses a ship if \", \" A player loses a ship if \"),\n    @(\"the high score is checked and all settings are reset\",\n      \"the high score is checked, and all settings are reset\"),\n    @(\"This file although very short allows\", \"This file, although very short allows\"),\n    @(\"As mentioned it also resets\", \"As mentioned, it also resets\"),\n    @(\"represents the ships location. The file has\",\n      \"represents the ships\" + [char]0x2019 + \" location. The file has\"),\n    @(\"one of which updates the ships position based on the players input,\",\n      \"one of which updates the ship\" + [char]0x2019 + \"s position based on the players input,\"),\n    @(\"speed up subsequent imports of the same module \",\n      \"speed up subsequent imports of the same module. \"),\n    @(\"Here we will talk about adjusting the speeds and sizes of bullets that we tried to get a challenging but fun user experience. While testing, during various steps of the build , we changed\",\n      \"Here we talked about and tested adjusting the speeds and sizes of bullets, we tried to get a challenging but fun user experience. While testing, during various steps of the build, we changed\"),\n    # \"complete\" is the unique remaining word of the \"Aug 5 - 10\" bullet; matched\n    # on its own (rather than together with the preceding sentence) so the\n    # existing grammar-check markers around it stay intact.\n    @(\"complete\",\n      \"complete, still unable to reach other members of group B\"),\n    @(\"By the end of the project it made sense\", \"By the end of the project, it made sense\"),\n    @(\"or any other real world tasks,\", \"or any other real-world tasks,\"),\n    @(\"We have also discussed about finding and implementing\",\n      \"We have also discussed finding and implementing\"),\n    @(\"maybe the bullets do no disappear\", \"maybe the bullets do not disappear\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n\nWrite-Output \"done\"\n"}
